$wb = $excel.ActiveWorkbook

# "Silver Rear_side": B9 5,263 -> 5,273
$ws1 = $wb.Worksheets.Item("Silver Rear_side")
$ws1.Range("B9").NumberFormat = "@"
$ws1.Range("B9").Value = "5,273"
$ws1.Range("B9").ClearFormats()

# "Silver Busbar front-side": B9 7,879 -> 7,895
$ws2 = $wb.Worksheets.Item("Silver Busbar front-side")
$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "7,895"
$ws2.Range("B9").ClearFormats()

# "Silver finger front-side": B9 7,929 -> 7,945
$ws3 = $wb.Worksheets.Item("Silver finger front-side")
$ws3.Range("B9").NumberFormat = "@"
$ws3.Range("B9").Value = "7,945"
$ws3.Range("B9").ClearFormats()

# "USD_CNY": B9 7.2647 -> 7.2597
$ws4 = $wb.Worksheets.Item("USD_CNY")
$ws4.Range("B9").NumberFormat = "@"
$ws4.Range("B9").Value = "7.2597"
$ws4.Range("B9").ClearFormats()
